$d = $word.ActiveDocument

# Update table values per 2021 Q2 Quarterly Report Final

$d.Content.Find.Execute("93 (30.4)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94 (30.6)", 2)

$d.Content.Find.Execute("94 (30.7)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "94 (30.6)", 2)

$d.Content.Find.Execute("91 (29.7)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "91 (29.6)", 2)

$d.Content.Find.Execute("306 (100.0)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "307 (100.0)", 2)
